$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the codelist labels in column A to use readable text instead of
# the slugified identifiers (per the diff on xl/sharedStrings.xml).
$ws.Range("A4").Value  = "10 a 19"
$ws.Range("A7").Value  = "40 o más"
$ws.Range("A8").Value  = "20 a 39"
$ws.Range("A9").Value  = "Sin determinar"
$ws.Range("A10").Value = "05 a 09"
